# Adds the new "ODI Bowling Extra" worksheet (sheetId=5) with its
# MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns, matching the
# commit "[ADDITIONAL SCRAPING] added scraping code for extra bowling
# attributes and excel sheets".

$wb = $excel.ActiveWorkbook

# --- Create the new sheet, positioned after the current last sheet ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "MATCH_CODE"
$ws.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$ws.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# --- Data rows -------------------------------------------------------------
# MATCH_CODE pulled from "ODI Bowling"; MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL
# left blank where the source had no value, kept as plain text otherwise
# (percentages are literal strings like "10.00%", not numeric cells).
$data = @(
    @{ Row = 2; A = "3997"; B = $null; C = $null },
    @{ Row = 3; A = "3999"; B = $null; C = $null },
    @{ Row = 4; A = "4002"; B = $null; C = $null },
    @{ Row = 5; A = "4029"; B = "0"; C = "10.00%" },
    @{ Row = 6; A = "4080"; B = "0"; C = "20.00%" },
    @{ Row = 7; A = "4083"; B = "0"; C = "20.00%" },
    @{ Row = 8; A = "4224"; B = "0"; C = "30.00%" },
    @{ Row = 9; A = "4226"; B = "0"; C = "20.00%" },
    @{ Row = 10; A = "4237"; B = "0"; C = $null },
    @{ Row = 11; A = "4247"; B = $null; C = $null },
    @{ Row = 12; A = "4261"; B = "0"; C = $null },
    @{ Row = 13; A = "4269"; B = "0"; C = $null },
    @{ Row = 14; A = "4272"; B = $null; C = $null },
    @{ Row = 15; A = "4303"; B = $null; C = $null },
    @{ Row = 16; A = "4339"; B = "2"; C = "30.00%" },
    @{ Row = 17; A = "4351"; B = "2"; C = "20.00%" },
    @{ Row = 18; A = "4529"; B = "0"; C = "20.00%" },
    @{ Row = 19; A = "4559"; B = "0"; C = $null },
    @{ Row = 20; A = "4619"; B = $null; C = $null },
    @{ Row = 21; A = "4620"; B = $null; C = $null }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = "'" + $item.A
    if ($item.B -ne $null) {
        $ws.Cells.Item($r, 2).Value = "'" + $item.B
    } else {
        # Keep the cell present-but-blank (as plain text), matching the
        # source sheet's empty MAIDEN_OVERS cells.
        $ws.Cells.Item($r, 2).Value = "'"
    }
    if ($item.C -ne $null) {
        $ws.Cells.Item($r, 3).Value = "'" + $item.C
    } else {
        $ws.Cells.Item($r, 3).Value = "'"
    }
}
